$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert the two new paragraphs right after "Contexte du projet :" and
#    before the pre-existing empty paragraph that follows it.
# ---------------------------------------------------------------------------

$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Contexte du projet*") {
        $targetIdx = $i
        break
    }
}

$emptyParaIdx = $targetIdx + 1

# Insert a bare paragraph (no inherited pPr/rsid) right before the empty one.
# This mimics Word's "Next style" behaviour when pressing Enter after a
# Titre paragraph (falls back to Normal / no explicit pPr).
$d.Paragraphs($emptyParaIdx).Range.InsertParagraphBefore()
$firstNewParaIdx = $emptyParaIdx
$emptyParaIdx = $emptyParaIdx + 1

# Build the first paragraph out of several runs (as in the source document)
# using InsertXML, which - unlike InsertAfter - keeps each run separate even
# though they all share identical (empty) run formatting.
$eacute = [char]0x00E9
$rsquo = [char]0x2019
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Notre objectif est de r' + $eacute + 'aliser un capteur capacitif en utilisant u</w:t></w:r><w:r><w:t>n l' + $eacute + 'gume</w:t></w:r><w:r><w:t xml:space="preserve">. Celui-ci </w:t></w:r><w:r><w:t xml:space="preserve">pourra </w:t></w:r><w:r><w:t>d' + $eacute + 'tecter</w:t></w:r><w:r><w:t xml:space="preserve"> si l' + $rsquo + 'on pose</w:t></w:r><w:r><w:t xml:space="preserve"> un doigt dessus, ou plusieurs afin de r' + $eacute + 'aliser diff' + $eacute + 'rentes fonctions. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs($firstNewParaIdx).Range.InsertXML($xml1)

# Insert the second bare paragraph before the (still) empty, pre-existing one.
$d.Paragraphs($emptyParaIdx).Range.InsertParagraphBefore()
$secondNewParaIdx = $emptyParaIdx
$emptyParaIdx = $emptyParaIdx + 1

$oelig = [char]0x0153
$egrave = [char]0x00E8
$d.Paragraphs($secondNewParaIdx).Range.InsertAfter("Pour cela, nous allons realiser differentes experiences qui nous aideront a mettre en " + $oelig + "uvre le syst" + $egrave + "me de detection")

# ---------------------------------------------------------------------------
# 2) Move the automatic "_GoBack" bookmark from the table cell (after "0.6")
#    to the end of the text we just typed in the second new paragraph.
# ---------------------------------------------------------------------------

$old = $d.Bookmarks("_GoBack")
$old.Delete()

# Insert a one-character placeholder right after the paragraph's text (but
# before its paragraph mark) so we can anchor a bookmark on a 1-char-wide
# range (the engine mishandles zero-length ranges passed to Bookmarks.Add).
# Deleting the placeholder afterwards collapses the bookmark to zero width
# exactly where we want it - matching real Word's behaviour.
$secPara = $d.Paragraphs($secondNewParaIdx)
$endRange = $secPara.Range
$endRange.MoveEnd(1, -1)
$endRange.Collapse(0)
$endRange.InsertAfter("X")

$secPara2 = $d.Paragraphs($secondNewParaIdx)
$bmRange = $secPara2.Range
$bmRange.MoveEnd(1, -1)
$bmRange.MoveStart(0, 1)
$bmRange.MoveStart(1, -1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$secPara3 = $d.Paragraphs($secondNewParaIdx)
$delRange = $secPara3.Range
$delRange.MoveEnd(1, -1)
$delRange.MoveStart(0, 1)
$delRange.MoveStart(1, -1)
$delRange.Delete()
